$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete entire column I ("Subscription ID"), shifting columns J:N left to I:M
$ws.Columns("I").Delete()

# The autofilter range still references the old last column (N) after the
# column delete, so drop it and re-apply over the new A1:M1 extent.
$ws.AutoFilterMode = $false
$ws.Range("A1:M1").AutoFilter()

# Keep the workbook-level _FilterDatabase defined name (used by the
# autofilter) in sync with the resized range as well.
for ($i = 1; $i -le $wb.Names.Count; $i++) {
  $n = $wb.Names.Item($i)
  if ($n.Name -like "*_FilterDatabase*") {
    $n.RefersTo = "=Data!`$A`$1:`$M`$1"
  }
}

# Match the author's final cursor position recorded in the saved file.
$ws.Range("J10").Select()
